$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Actividades durante la pasantia")

# Fix the array formula in A2 (30/10/2021 -> 31/10/2021 typo correction)
$ws.Range("A2").FormulaArray = "=A2:A2931/10/2021"

# Copy formatting of the last existing data row (105) down into the new rows
# so the new rows inherit the same styles (date format, wrap text, borders, etc.)
$ws.Range("A105:C105").Copy()
$ws.Range("A106:C111").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the shared-string text cells in the same order the original
# author entered them, so new entries land at the same shared-string indices.
$ws.Cells.Item(107,2).Value = "Discusion articulo"
$ws.Cells.Item(106,2).Value = "Busqueda de articulos"
$ws.Cells.Item(108,2).Value = "Discusion articulo"
$ws.Cells.Item(109,2).Value = "Busqueda de articulos ASC, reunion vero"
$ws.Cells.Item(110,2).Value = "Revision huber con vero, articulos y preguntas ASC para los retos explora"
$ws.Cells.Item(111,2).Value = "Obtener dataframes por columnas de SL, coherencia"

# New activity rows - dates and hours
$ws.Cells.Item(106,1).Value = 44843
$ws.Cells.Item(106,3).Value = 3

$ws.Cells.Item(107,1).Value = 44844
$ws.Cells.Item(107,3).Value = 8

$ws.Cells.Item(108,1).Value = 44845
$ws.Cells.Item(108,3).Value = 2

$ws.Cells.Item(109,1).Value = 44846
$ws.Cells.Item(109,3).Value = 2

$ws.Cells.Item(110,1).Value = 44849
$ws.Cells.Item(110,3).Value = 6
$ws.Rows.Item(110).RowHeight = 27.6

$ws.Cells.Item(111,1).Value = 44852
$ws.Cells.Item(111,3).Value = 8

# Update the view: scroll down and select the new last cell, like the author did
$ws.Range("C111").Select()
$excel.ActiveWindow.ScrollRow = 101
$excel.ActiveWindow.ScrollColumn = 1
